# Saldo_guide.xlsx update
# - Advance the "Dt. Referencia" (column G) for every data row from 45406 (2024-04-24)
#   to 45407 (2024-04-25).
# - A handful of accounts had their Saldo Previsto / Vl. Total (columns D and H)
#   recalculated for the new reference date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference date for every data row (rows 2 through 310).
for ($r = 2; $r -le 310; $r++) {
    $ws.Cells.Item($r, 7).Value = 45407
}

# Rows whose Saldo Previsto (D) / Vl. Total (H) values changed for the new date.
$ws.Cells.Item(12, 4).Value = 26589.35
$ws.Cells.Item(12, 8).Value = 26589.35

$ws.Cells.Item(71, 4).Value = 4986.86
$ws.Cells.Item(71, 8).Value = 4986.86

$ws.Cells.Item(78, 4).Value = 0
$ws.Cells.Item(78, 8).Value = 0

$ws.Cells.Item(89, 4).Value = 0
$ws.Cells.Item(89, 8).Value = 0

$ws.Cells.Item(184, 4).Value = 37189.27
$ws.Cells.Item(184, 8).Value = 37189.27

Write-Host "Updated reference dates and balances."
